$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15-22 down to 16-23
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with data
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = 44512
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100101
$ws.Cells.Item(15, 8).Value = "Berries"
$ws.Cells.Item(15, 9).Value = 100101001
$ws.Cells.Item(15, 10).Value = "Arándano (blue)"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 30
$ws.Cells.Item(15, 14).Value = 6000
$ws.Cells.Item(15, 15).Value = 6000
$ws.Cells.Item(15, 16).Value = 6000
$ws.Cells.Item(15, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Linares"
$ws.Cells.Item(15, 19).Value = 3000
$ws.Cells.Item(15, 20).Value = 2

# Insert a second new row at row 21 (after the first insert shifted rows down),
# shifting rows 21-23 down to 22-24
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with data
$ws.Cells.Item(21, 1).Value = 5
$ws.Cells.Item(21, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(21, 3).Value = "Maule"
$ws.Cells.Item(21, 4).Value = 44511
$ws.Cells.Item(21, 5).Value = 7
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100101
$ws.Cells.Item(21, 8).Value = "Berries"
$ws.Cells.Item(21, 9).Value = 100101001
$ws.Cells.Item(21, 10).Value = "Arándano (blue)"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 50
$ws.Cells.Item(21, 14).Value = 6400
$ws.Cells.Item(21, 15).Value = 6400
$ws.Cells.Item(21, 16).Value = 6400
$ws.Cells.Item(21, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(21, 18).Value = "Provincia de Linares"
$ws.Cells.Item(21, 19).Value = 3200
$ws.Cells.Item(21, 20).Value = 2

# Ensure the date cells carry the date number format, matching column D elsewhere
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
